$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 270 (shifts IAD..YHZ down by one row)
$ws.Rows.Item(270).Insert()

# Copy the (same) formatting from the colo-code cell below onto the new row's
# colo-code cell, so it keeps the bold/bordered "colo" column style.
$ws.Range("A271").Copy()
$ws.Range("A270").PasteSpecial(-4122)

# Populate the newly inserted row 270 with the MLG (Malang, Indonesia) colo data
$ws.Range("A270").Value = "MLG"
$ws.Range("B270").Value = "Malang, Indonesia"
$ws.Range("C270").Value = "Asia Pacific"
$ws.Range("D270").Value = "Malang"
$ws.Range("E270").Value = "Indonesia"
$ws.Range("F270").Value = "ID"
$ws.Range("G270").Value = -8.100346999999999
$ws.Range("H270").Value = 112.186641
